# Update "想去人数" (F column) figures across the 展览 / 演出 / 全部类型
# sheets to the freshly scraped counts (gh-pages output regenerated at
# 456a3b4). 本地生活 has no data rows, so it needs no changes.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExpo.Range("F3").Value  = 2172
$wsExpo.Range("F4").Value  = 50
$wsExpo.Range("F5").Value  = 11496
$wsExpo.Range("F6").Value  = 204
$wsExpo.Range("F7").Value  = 320
$wsExpo.Range("F8").Value  = 221
$wsExpo.Range("F9").Value  = 11442
$wsExpo.Range("F10").Value = 468
$wsExpo.Range("F11").Value = 1159
$wsExpo.Range("F12").Value = 78
$wsExpo.Range("F13").Value = 1751
$wsExpo.Range("F14").Value = 5688
$wsExpo.Range("F16").Value = 3488

# 演出 (sheet2)
$wsShow.Range("F2").Value = 569
$wsShow.Range("F4").Value = 2

# 全部类型 (sheet4) - combined view of 展览 + 演出 rows
$wsAll.Range("F3").Value  = 2172
$wsAll.Range("F4").Value  = 569
$wsAll.Range("F5").Value  = 50
$wsAll.Range("F7").Value  = 11496
$wsAll.Range("F8").Value  = 204
$wsAll.Range("F9").Value  = 320
$wsAll.Range("F10").Value = 221
$wsAll.Range("F11").Value = 11442
$wsAll.Range("F12").Value = 468
$wsAll.Range("F13").Value = 1159
$wsAll.Range("F14").Value = 78
$wsAll.Range("F15").Value = 1751
$wsAll.Range("F16").Value = 2
$wsAll.Range("F17").Value = 5688
$wsAll.Range("F19").Value = 3488
